# query options: replace: works, start on postprocess:
#
# SOURCE_A.xlsx has two sheets:
#   "A.DATETIME" (sheet1) - view/selection only changes
#   "A.1"        (sheet2) - gains a second column ("FIELD_A") and its
#                 A-column values flip from the literal string "A1.1"
#                 to real row numbers (1, 2) that the new B column
#                 ("MATCH A1") keys off of.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("A.DATETIME")
$ws2 = $wb.Worksheets.Item("A.1")

# --- "A.1" sheet: new FIELD_A column + numeric id column -------------
$ws2.Range("B1").Value = "FIELD_A"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "MATCH A1"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "MATCH A1"

# --- view / selection state -------------------------------------------
# "A.DATETIME" keeps its data but is no longer the active tab; its
# selection moves from D2 to E1.
$ws1.Range("E1").Select()

# "A.1" becomes the active sheet/tab, selection moves from A4 to A2.
$ws2.Activate()
$ws2.Range("A2").Select()
